$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("YDS")
$ws.Range("B2").Value = "5 -2 3 4 1 2 -1 13 1 3 6 4 7 2 2 5 10 19 1 5 5 9 12 -1 7 1 2 8 2 4 0 3 0 5 2 3 2 3 1 8 -1 -3 4 9 5 1 5 10 3 3 5 2 4 3 2 6 4 -2 6 27 1 -1 3 11 11 4 3 2 8 1 1 5 2 7 -1 2 2 2 7 1 0 4 3 9 4 11 5 7 0 1 4 0 1 8 1 2 4 0 7 10 -1 9 4 0 6 5 7 11 1 -1 4 9 8 2 3 -2 10 4 4 1 1 0 9 8 5 2 3 0 2 4 6 6 14 4 3 5 1 7 6 3 4 6 2 4 5 18 4 2 1 -2 6 4 2 9 0 -1 -1 9 0 0 5 1 0 0 10 7 21 3 46 3 -3 -5 2 2 1 -1 -1 -1 1 6 8 0 -4 6 7 4 9 4 5 -1 2 1 2 0 3 3 7 4 3 0 21 2 2 4 0 1 2 5 0 2 9 2 6 4 14 4 5 1 3 -2 -3 1 -2 10 -1 11 2 5 4 8 -2 18 5 2 6 4 -5 2 1 19 4 0 1 3 7 20 1 0 1 25 3 7 2 3 6 3 8 3 1 6 4 5 9 8 9 0 9 5 1 8 9 4 -1 0 1 32 1 1 2 3 13 2 3 -2 -1 4 0 2 2 7 1 -5 6 8 6 2 6 17 2 5 3 4 0 10 3 3 8 -3 5 10 -2 8 3 3 3 2 4 6 4 0 2 7 1 6 1 2 4 6 2 0 17 1 0 4 3 6 0 2 4 5 5 6 5 1 -1 5 8 4 3 5 22 2 12 1 8 2 4 1 1 3 8 -2 19 3 1 1 9 1 1 1 3 12 2 2 2 8 2 2 2 13 3 6 6 3 1 4 3 3 1 3 4 5 4 1 0 0 0 3 4 3 5 2 5 7 2 4 -1 9 1 10 5 4 2 2 2 3 0 -2 1 1 1 2 5 -3 0 1 2 6 10 7 16 1 13 -1 4 2 -2 3 3 8 1 23 12 3 1 -1 -2 2 5 2 0 0 5 3 4 7 1 11 4 2 7 13 2 -3"
$ws.Range("C2").Value = "1 3 1 2 17 2 7 -3 4 1 13 0 8 -1 6 5 2 1 0 -1 0 0 2 2 3 7 8 6 11 14 5 0 -1 0 6 6 5 11 7 2 0 -3 6 1 1 4 4 -1 1 10 6 4 4 2 4 -3 5 1 1 3 20 0 2 2 1 2 1 11 3 6 2 4 11 2 11 12 6 1 3 5 5 20 1 3 0 7 2 3 0 8 0 2 4 1 4 2 14 11 3 0 2 1 3 8 1 5 7 17 -1 5 4 1 1 57 0 -1 -4 -1 3 1 0 3 8 1 6 2 1 1 -5 1 0 9 5 -2 1 2 10 -5 4 -2 -1 12 3 10 3 7 1 16 2 1 7 7 2 10 7 2 5 6 0 6 8 6 2 5 9 -1 1 8 4 4 3 2 5 5 18 2 -1 1 4 2 7 1 11 22 0 1 2 3 0 1 70 7 1 -2 13 10 2 3 2 -1 -1 -4 6 4 3 3 6 4 5 1 7 4 1 4 3 2 4 5 0 18 -4 7 2 3 1 2 2 4 5 4 3 5 0 4 3 -1 13 3 2 0 -1 6 7 3 9 5 0 2 5 4 1 15 1 5 6 1 3 -4 7 2 -1 0 1 4 4 9 12 2 27 1 7 2 3 -3 2 5 3 3 1 8 6 2 2 4 2 4 4 3 5 -1 5 2 2 1 3 2 6 7 2 2 5 6 14 4 5 7 1 0 -1 7 7 13 1 -2 10 1 2 8 2 6 1 0 7 -2 3 -5 3 4 2 2 -7 1 1 2 23 -1 6 2 1 8 10 6 8 1 14 3 17 1 19 8 0 5 8 11 5 3 1 2 2 14 35 2 20 -1 0 8 7 5 2 -1 8 3 -1 2 4 3 4 6 3 8 9 3 9 45 5 5 1 -1 9 10 0 3 0 8 1 1 8 11 7 7 1 3 7 14 4 6 4 9 4 7 4 1 9 14 4 4 1"
$ws.Range("B3").Value = "5 3 7 12 17 13 2 14 2 9 50 4 28 5 14 15 18 5 6 2 7 12 10 3 12 1 2 14 8 7 14 14 22 5 10 15 42 7 12 15 7 17 7 19 34 14 8 14 5 9 4 7 33 -2 15 18 0 10 21 0 44 22 25 12 1 11 8 4 8 12 7 31 12 15 9 7 25 2 7 9 6 7 2 10 4 7 5 70 28 14 5 8 3 10 2 7 19 9 8 8 5 21 14 24 5 5 1 9 4 3 7 9 34 8 7 23 11 53 14 2 7 11 23 6 5 9 55 7 4 6 13 26 7 17 21 32 7 82 25 2 11 11 27 6 3 13 0 3 9 16 4 26 2 46 9 21 8 54 8 10 11 6 21 10 19 11 11 9 20 15 19 16 9 9 5 11 16 6 12 5 5 12 8 11 9 3 3 4 11 11 11 9 11 6 3 6 3 2 7 14 9 14 3 7 6 4 17 7 3 7 9 9 12 6 5 9 5 16 4 12 32 11 16 7 15 16 6 0 -4 0 -7 16 12 19 5 12 21 3 9 7 32 16 29 5 7 0 5 9 20 12 6 9 11 12 16 14 25 14 10 14 20 0 10 11 16 9 12 -4 11 12 19 3 5 27 4 13 5 12 22 29 32 26 23 6 2 6 17 2 2 9 8 4 19 3 56 4 4 15 10 18 20 3 11 16 5 17 3 68 20 17 9 11 15 3 52 1 19 3 18 9 3 2 5 14 13 0 12 4 22 3 10 20 4 52 11 4 3 10 72 8 8 9 17 18 6 7 6 4 10 5 69 22 6 4 39 10 5 2 5 9 2 35 30 11 17 9 21 7 -1 28 5 5 29 5 10 5 9 7 9 5 7 3 19 11 10 6 7 21 8 4 57 7 19 7 7 6 15 4 22 1 5 10 11 32 7 12 11 4 7 21 8 13 7 3 19 2 2 7 3 5 9 6 8 18 41 44 4 22 18 2 17 2 2 16 5 0 9 8"
$ws.Range("C3").Value = "9 6 12 11 9 -1 2 9 4 17 11 25 5 5 6 20 0 9 10 10 5 6 9 8 34 10 7 27 0 6 14 7 3 5 12 10 17 -5 3 3 11 2 36 10 5 11 10 13 21 5 0 5 3 0 6 3 7 3 16 12 6 6 3 18 4 6 15 21 16 4 24 0 6 8 3 5 8 7 5 8 11 17 5 29 9 3 8 8 8 17 7 12 10 7 7 52 7 4 27 5 20 -2 7 11 5 15 17 -4 34 14 0 12 13 24 -1 3 15 24 5 22 5 7 0 8 8 2 2 59 20 20 15 0 7 10 4 6 5 3 5 -2 5 4 33 2 9 5 3 15 7 4 3 7 10 3 7 8 10 14 15 24 6 23 19 20 18 7 1 35 39 25 3 13 13 -1 10 12 7 11 9 30 8 7 6 3 15 3 6 7 12 6 20 6 5 2 26 11 12 15 9 4 8 13 14 15 9 8 9 10 24 14 23 12 6 8 8 19 13 3 23 3 11 60 15 6 5 22 7 -2 21 26 18 -2 25 6 8 6 21 4 1 9 2 25 31 19 7 8 20 9 8 8 11 41 0 11 7 11 2 15 3 10 12 2 30 9 16 13 9 6 8 15 6 5 10 11 -1 10 -1 41 5 10 47 6 7 44 23 7 4 7 12 7 12 9 7 6 -3 1 18 33 2 7 5 4 2 10 21 6 1 4 3 22 20 16 10 13 8 9 5 17 19 8 25 21 3 9 12 5 8 7 28 -1 10 3 11 -1 8 9 11 8 25 2 24 24 5 11 17 5 15 6 4 9 -1 16 12 10 -1 18 28 11 11 4 7 12 17 25 18 14 4 6 7 14 5 8 6 1 29 7 17 24 3 3 53 7 3 2 4 27 7 8 5 11 6 5 6 4 8 1 4 17 6 8 9 13 7 6 7 19 4 11 14 12 6 9 16 7 5 16 3 3 13 5 13 7 15 23 10 11 13 3 41 3 13 8 40 5 33 20 6 3 16 5 14 8 10 6 10 44 11 7 5 19 5 12 3 33 14 9 10 0 6 4 7 5 5 8 18 2"

$ws = $wb.Worksheets.Item("OFF")
$ws.Range("C2").Value = 245
$ws.Range("D2").Value = 11
$ws.Range("F2").Value = 61
$ws.Range("G2").Value = 69
$ws.Range("I2").Value = 7
$ws.Range("J2").Value = 31
$ws.Range("N2").Value = 26
$ws.Range("B3").Value = 10
$ws.Range("C3").Value = 209
$ws.Range("D3").Value = 7
$ws.Range("E3").Value = 47
$ws.Range("F3").Value = 134
$ws.Range("G3").Value = 44
$ws.Range("H3").Value = 47
$ws.Range("I3").Value = 59
$ws.Range("J3").Value = 72
$ws.Range("L3").Value = 275
$ws.Range("M3").Value = 184
$ws.Range("Q3").Value = 532

$ws = $wb.Worksheets.Item("DEF")
$ws.Range("C2").Value = 208
$ws.Range("D2").Value = 11
$ws.Range("F2").Value = 56
$ws.Range("G2").Value = 62
$ws.Range("J2").Value = 31
$ws.Range("N2").Value = 25
$ws.Range("C3").Value = 254
$ws.Range("D3").Value = 4
$ws.Range("F3").Value = 137
$ws.Range("G3").Value = 59
$ws.Range("I3").Value = 72
$ws.Range("J3").Value = 70
$ws.Range("L3").Value = 339
$ws.Range("M3").Value = 219
$ws.Range("Q3").Value = 586

$ws = $wb.Worksheets.Item("ST")
$ws.Range("B2").Value = 108
$ws.Range("D2").Value = 68
$ws.Range("F2").Value = 52
$ws.Range("G2").Value = 50
$ws.Range("H2").Value = 5
$ws.Range("I2").Value = 2
$ws.Range("J2").Value = 59
$ws.Range("K2").Value = 57
$ws.Range("N2").Value = 24
$ws.Range("O2").Value = 19
$ws.Range("B3").Value = 67
$ws.Range("D3").Value = "59 47 37 43 42 43 61 36 39 48 37 54 49 53 46 56 40 46 41 38 42 25 58 51 54 39 60 58 59 47 45 38 55 37 44 41 47 58 41 45 53 50 46 46 38 46 36 41 50 57 57 41 61 45 56 51 40 51 27 33 48 39 43 31 42 54 31 36"
$ws.Range("B4").Value = "70 55 63 69 56 62 67 62 62 66 64 67 64 39 62 68 62 66 59 67 68 68 65 55 63 67 66 56 64 63 64 64 67 45 63 61 61 64 65 62 66"
$ws.Range("D4").Value = "0 0 0 0 0 0 12 0 0 14 0 0 0 0 0 0 0 1 0 0 8 0 11 13 14 0 16 21 0 0 0 0 16 0 -1 5 12 11 0 3 0 15 3 0 0 0 4 0 0 0 0 0 0 0 12 18 0 29 3 0 0 5 0 0 0 0 0 0"
$ws.Range("B5").Value = "21 10 22 27 16 20 23 28 18 21 24 35 24 11 12 25 15 24 24 25 48 24 27 16 24 23 22 15 24 0 25 26 30 12 0 35 16 22 26 16 22"
$ws.Range("D5").Value = "0 0 12 8 0 1 0 12 0 0 5 0 5 6 10 0 0 3 0 14 12 0 0 0 0 11 0 11 0"
$ws.Range("B6").Value = "11 27 18 23 30 19 24 44 22 28 26 18 14 21 17 22 25 32 12 17 10 25 18 18 23 25 20 0 14 27 23 25 31 32 19 26"

$ws = $wb.Worksheets.Item("TURNS")
$ws.Range("B3").Value = 9
$ws.Range("C3").Value = 12
$ws.Range("E3").Value = 12

$ws = $wb.Worksheets.Item("PEN")
$ws.Range("B2").Value = 16
$ws.Range("D2").Value = 7
$ws.Range("D4").Value = 9
